$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 322
$ws.Range("I2").Value = 845
$ws.Range("J2").Value = 3683
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 952
$ws.Range("M2").Value = 57
$ws.Range("N2").Value = 617
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 17
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 51
$ws.Range("S2").Value = 447
$ws.Range("T2").Value = 618
$ws.Range("U2").Value = 47
$ws.Range("V2").Value = 5540
$ws.Range("X2").Value = 5705
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 92
$ws.Range("AA2").Value = 37
